$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P1").Value = "ExpectedRecordDuration"
$ws.Range("Q1").Value = "OmicronFile"
$ws.Range("Q2").Value = "'CAM_763.seq"

$ws.Range("P1").Select()
